$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 197.5433703333333
$ws.Range("H2").Value = 592.6301109999999
$ws.Range("I2").Value = 0.3388703761585983
$ws.Range("J2").Value = 0.3388703761585982
$ws.Range("M2").Value = 1.223926333333333
$ws.Range("N2").Value = 3.671779
$ws.Range("Q2").Value = 241.7785329263855
$ws.Range("R2").Value = 2176.006796337469
$ws.Range("S2").Value = 0.3388703761585983
$ws.Range("T2").Value = 0.3388703761585982

$ws.Range("I3").Value = 0.1369374790620155
$ws.Range("J3").Value = 0.1369374790620154
$ws.Range("M3").Value = 1.223926333333333
$ws.Range("N3").Value = 3.671779
$ws.Range("Q3").Value = 97.7026766563869
$ws.Range("R3").Value = 879.3240899074821
$ws.Range("S3").Value = 0.1369374790620155
$ws.Range("T3").Value = 0.1369374790620154

$ws.Range("G4").Value = 148.824417
$ws.Range("H4").Value = 446.473251
$ws.Range("I4").Value = 0.2552967790580629
$ws.Range("J4").Value = 0.2552967790580629
$ws.Range("M4").Value = 1.223926333333333
$ws.Range("N4").Value = 3.671779
$ws.Range("Q4").Value = 182.150123009281
$ws.Range("R4").Value = 1639.351107083529
$ws.Range("S4").Value = 0.2552967790580629
$ws.Range("T4").Value = 0.2552967790580629

$ws.Range("G5").Value = 35.426853
$ws.Range("H5").Value = 106.280559
$ws.Range("I5").Value = 0.06077202683121193
$ws.Range("J5").Value = 0.06077202683121192
$ws.Range("M5").Value = 1.223926333333333
$ws.Range("N5").Value = 3.671779
$ws.Range("Q5").Value = 43.35985829382901
$ws.Range("R5").Value = 390.238724644461
$ws.Range("S5").Value = 0.06077202683121193
$ws.Range("T5").Value = 0.06077202683121192

$ws.Range("G6").Value = 121.3248153333333
$ws.Range("H6").Value = 363.974446
$ws.Range("I6").Value = 0.2081233388901116
$ws.Range("J6").Value = 0.2081233388901115
$ws.Range("M6").Value = 1.223926333333333
$ws.Range("N6").Value = 3.671779
$ws.Range("Q6").Value = 148.4926363732705
$ws.Range("R6").Value = 1336.433727359434
$ws.Range("S6").Value = 0.2081233388901116
$ws.Range("T6").Value = 0.2081233388901115
